# Rename "Test samples" sheet to "Protocol Testing"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test samples")
$ws.Name = "Protocol Testing"

# Add new gel-notes rows to the bottom of the sheet (rows 36-37 and 39-40)
$ws.Range("B36").Value = "GeneRuler 1kb plus DNA LAdder on 2% gel: 100V, 400amp, 1 hour"
$ws.Range("A36").Value = "Gel 1:"

$ws.Range("A39").Value = "Gel 2:"

$ws.Range("A37").Value = "Above 30 samples from all 3 temps"

$ws.Range("A40").Value = "57C set, with bright band samples diluted and not diluted to see if there are 2 bands "

$ws.Range("B39").Value = "GeneRuler 100 bp DNA Ladder on 2% gel: 75V, 400amp, 1 hour"

# Update the active selection on this sheet to match the authored state
$ws.Range("E16").Select()
